$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5380440000000001
$ws.Range("H2").Value = 1.614132
$ws.Range("I2").Value = 0.9003438764610565
$ws.Range("J2").Value = 0.9003438764610566
$ws.Range("M2").Value = 0.4102596666666667
$ws.Range("N2").Value = 1.230779
$ws.Range("O2").Value = 0.003499619873322347
$ws.Range("P2").Value = 0.003499619873322347
$ws.Range("Q2").Value = 0.2207377520920001
$ws.Range("R2").Value = 1.986639768828
$ws.Range("S2").Value = 0.003150861322887194
$ws.Range("T2").Value = 0.003150861322887194
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5380440000000001
$ws.Range("H3").Value = 1.614132
$ws.Range("I3").Value = 0.9003438764610565
$ws.Range("J3").Value = 0.9003438764610566
$ws.Range("O3").Value = 0.8692174743460166
$ws.Range("P3").Value = 0.8692174743460165
$ws.Range("Q3").Value = 54.82570059361201
$ws.Range("R3").Value = 493.431305342508
$ws.Range("S3").Value = 0.7825946303403816
$ws.Range("T3").Value = 0.7825946303403816
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5380440000000001
$ws.Range("H4").Value = 1.614132
$ws.Range("I4").Value = 0.9003438764610565
$ws.Range("J4").Value = 0.9003438764610566
$ws.Range("N4").Value = 44.764041
$ws.Range("O4").Value = 0.1272829057806611
$ws.Range("P4").Value = 0.1272829057806611
$ws.Range("Q4").Value = 8.028341225268001
$ws.Range("R4").Value = 72.255071027412
$ws.Range("S4").Value = 0.1145983847977879
$ws.Range("T4").Value = 0.1145983847977879
$ws.Range("G5").Value = 0.05955433333333333
$ws.Range("I5").Value = 0.09965612353894335
$ws.Range("J5").Value = 0.09965612353894336
$ws.Range("M5").Value = 0.4102596666666667
$ws.Range("N5").Value = 1.230779
$ws.Range("O5").Value = 0.003499619873322347
$ws.Range("P5").Value = 0.003499619873322347
$ws.Range("Q5").Value = 0.02443274094188889
$ws.Range("R5").Value = 0.219894668477
$ws.Range("S5").Value = 0.0003487585504351531
$ws.Range("T5").Value = 0.0003487585504351531
$ws.Range("G6").Value = 0.05955433333333333
$ws.Range("I6").Value = 0.09965612353894335
$ws.Range("J6").Value = 0.09965612353894336
$ws.Range("O6").Value = 0.8692174743460166
$ws.Range("P6").Value = 0.8692174743460165
$ws.Range("Q6").Value = 6.068477760899666
$ws.Range("R6").Value = 54.616299848097
$ws.Range("S6").Value = 0.08662284400563496
$ws.Range("T6").Value = 0.08662284400563496
$ws.Range("G7").Value = 0.05955433333333333
$ws.Range("I7").Value = 0.09965612353894335
$ws.Range("J7").Value = 0.09965612353894336
$ws.Range("N7").Value = 44.764041
$ws.Range("O7").Value = 0.1272829057806611
$ws.Range("P7").Value = 0.1272829057806611
$ws.Range("Q7").Value = 0.8886308730203332
$ws.Range("R7").Value = 7.997677857183
$ws.Range("S7").Value = 0.01268452098287325
$ws.Range("T7").Value = 0.01268452098287326
